$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs_for_tollcalib")
$ws.Range("A16:A164").Value = "RTP2021"
